# Updating filtered feeds from workflow
# Adds a new feed row (row 22) to the "Filtered Feeds" sheet:
#   A22: link to the new article (fiercebiotech)
#   B22: keyword "KRAS"
#   C22: title markup (<a href="..." hreflang="en">...</a>)
# and wires up a real hyperlink (like the other rows) pointing at the link.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newLink  = "https://www.fiercebiotech.com/biotech/bridgebios-kras-focused-cancer-spinout-backs-spac-path-nasdaq"
$newKeyword = "KRAS"
$newTitle = '<a href="https://www.fiercebiotech.com/biotech/bridgebios-kras-focused-cancer-spinout-backs-spac-path-nasdaq" hreflang="en">BridgeBio''s KRAS-focused cancer spinout backs SPAC as path to Nasdaq</a>'

$row = 22

$ws.Range("A" + $row).Value = $newLink
$ws.Range("B" + $row).Value = $newKeyword
$ws.Range("C" + $row).Value = $newTitle

# Wire up the hyperlink relationship (mirrors the other rows in column A),
# then re-apply the standard "Hyperlink" cell style used by the rest of the
# column so it matches the existing rows exactly.
$ws.Hyperlinks.Add($ws.Range("A" + $row), $newLink)
$ws.Range("A" + $row).Style = "Hyperlink"
